# Clean System Level v0.1
#
# Applies:
#  - baseline!E2:E5 values updated
#  - emission!B2:AA2 all set to a flat 10,000,000 (formulas in C2:G2 replaced
#    by plain values; H2:AA2 bumped from 200000 to 10000000)
#  - active-sheet / selection bookkeeping moves from "emission" to "baseline"

$wb = $excel.ActiveWorkbook

# --- baseline sheet -------------------------------------------------------
$baseline = $wb.Worksheets.Item("baseline")

$baseline.Range("E2").Value = 6000
$baseline.Range("E3").Value = 5000
$baseline.Range("E4").Value = 2000
$baseline.Range("E5").Value = 2000

# --- emission sheet ---------------------------------------------------------
$emission = $wb.Worksheets.Item("emission")

$emission.Range("B2:AA2").Value = 10000000

# --- selection / active sheet bookkeeping ----------------------------------
$emission.Activate()
$emission.Range("B2:AA2").Select()

$baseline.Activate()
$baseline.Range("E6").Select()
